$d = $word.ActiveDocument

function Replace-WithXml($searchText, $innerXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    $target = $d.Range($rng.Start, $rng.End)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
}

# 1. "If using an Xbox One devkit..." paragraph: split into 3 runs with a
#    gramStart/gramEnd proofErr pair bracketing "Gaming.Xbox.XboxOne.x".
Replace-WithXml `
    "If using an Xbox One devkit, set the active solution platform to Gaming.Xbox.XboxOne.x64." `
    ('<w:r><w:t xml:space="preserve">If using an Xbox One devkit, set the active solution platform to </w:t></w:r>' + `
     '<w:proofErr w:type="gramStart"/>' + `
     '<w:r><w:t>Gaming.Xbox.XboxOne.x</w:t></w:r>' + `
     '<w:proofErr w:type="gramEnd"/>' + `
     '<w:r><w:t>64.</w:t></w:r>')

# 2. "If using Project Scarlett..." paragraph: text changes to reference an
#    Xbox Series X|S devkit, split into 5 runs with a gramStart/gramEnd
#    proofErr pair bracketing "Gaming.Xbox.Scarlett.x".
Replace-WithXml `
    "If using Project Scarlett, set the active solution platform to Gaming.Xbox.Scarlett.x64." `
    ('<w:r><w:t xml:space="preserve">If using </w:t></w:r>' + `
     '<w:r><w:t>an Xbox Series X|S devkit</w:t></w:r>' + `
     '<w:r><w:t xml:space="preserve">, set the active solution platform to </w:t></w:r>' + `
     '<w:proofErr w:type="gramStart"/>' + `
     '<w:r><w:t>Gaming.Xbox.Scarlett.x</w:t></w:r>' + `
     '<w:proofErr w:type="gramEnd"/>' + `
     '<w:r><w:t>64.</w:t></w:r>')

# 3. "...a XUID to indicate who's friends..." run: split into 3 runs with a
#    gramStart/gramEnd proofErr pair bracketing "who's".
Replace-WithXml `
    " value passed into the query and including a XUID to indicate who’s friends should be included in the result list.  " `
    ('<w:r><w:t xml:space="preserve"> value passed into the query and including a XUID to indicate </w:t></w:r>' + `
     '<w:proofErr w:type="gramStart"/>' + `
     '<w:r><w:t>who’s</w:t></w:r>' + `
     '<w:proofErr w:type="gramEnd"/>' + `
     '<w:r><w:t xml:space="preserve"> friends should be included in the result list.  </w:t></w:r>')

Write-Output "done"
